$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 3).Value = 316626
$ws.Cells.Item(2, 4).Value = 403549255
$ws.Cells.Item(10, 3).Value = 116128
$ws.Cells.Item(10, 4).Value = 170165077
$ws.Cells.Item(12, 3).Value = 58664
$ws.Cells.Item(12, 4).Value = 84664847
$ws.Cells.Item(16, 3).Value = 3973
$ws.Cells.Item(16, 4).Value = 5638519
$ws.Cells.Item(20, 3).Value = 6504
$ws.Cells.Item(20, 4).Value = 9071285
$ws.Cells.Item(22, 3).Value = 76510
$ws.Cells.Item(22, 4).Value = 95459790
$ws.Cells.Item(28, 3).Value = 32229
$ws.Cells.Item(28, 4).Value = 47183637
$ws.Cells.Item(30, 3).Value = 11365
$ws.Cells.Item(30, 4).Value = 16346496
$ws.Cells.Item(35, 3).Value = 1784
$ws.Cells.Item(35, 4).Value = 2517833
$ws.Cells.Item(36, 3).Value = 96120
$ws.Cells.Item(36, 4).Value = 121044469
$ws.Cells.Item(44, 3).Value = 44080
$ws.Cells.Item(44, 4).Value = 64601423
$ws.Cells.Item(46, 3).Value = 9043
$ws.Cells.Item(46, 4).Value = 12978599
$ws.Cells.Item(48, 3).Value = 1396
$ws.Cells.Item(48, 4).Value = 1938609
$ws.Cells.Item(51, 3).Value = 2259
$ws.Cells.Item(51, 4).Value = 3152001
$ws.Cells.Item(52, 3).Value = 68255
$ws.Cells.Item(52, 4).Value = 85652120
$ws.Cells.Item(58, 3).Value = 27924
$ws.Cells.Item(58, 4).Value = 40953031
$ws.Cells.Item(61, 3).Value = 10963
$ws.Cells.Item(61, 4).Value = 15851456
$ws.Cells.Item(63, 3).Value = 1352
$ws.Cells.Item(63, 4).Value = 1890289
$ws.Cells.Item(67, 3).Value = 1443
$ws.Cells.Item(67, 4).Value = 2020565
$ws.Cells.Item(69, 3).Value = 20265
$ws.Cells.Item(69, 4).Value = 26541812
$ws.Cells.Item(73, 3).Value = 7521
$ws.Cells.Item(73, 4).Value = 11011030
$ws.Cells.Item(75, 3).Value = 5056
$ws.Cells.Item(75, 4).Value = 7341264
$ws.Cells.Item(78, 3).Value = 139074
$ws.Cells.Item(78, 4).Value = 173458789
$ws.Cells.Item(84, 3).Value = 63029
$ws.Cells.Item(84, 4).Value = 92382448
$ws.Cells.Item(87, 3).Value = 29386
$ws.Cells.Item(87, 4).Value = 42507898
$ws.Cells.Item(89, 3).Value = 2719
$ws.Cells.Item(89, 4).Value = 3915742
$ws.Cells.Item(90, 3).Value = 2768
$ws.Cells.Item(90, 4).Value = 3910279
$ws.Cells.Item(91, 3).Value = 32245
$ws.Cells.Item(91, 4).Value = 43682165
$ws.Cells.Item(95, 3).Value = 7820
$ws.Cells.Item(95, 4).Value = 11498881
$ws.Cells.Item(97, 3).Value = 7130
$ws.Cells.Item(97, 4).Value = 10336819
$ws.Cells.Item(99, 3).Value = 526
$ws.Cells.Item(99, 4).Value = 748405
$ws.Cells.Item(101, 3).Value = 8855
$ws.Cells.Item(101, 4).Value = 12290720
$ws.Cells.Item(103, 3).Value = 2227
$ws.Cells.Item(103, 4).Value = 3280802
$ws.Cells.Item(105, 3).Value = 2996
$ws.Cells.Item(105, 4).Value = 4376451
$ws.Cells.Item(109, 3).Value = 139528
$ws.Cells.Item(109, 4).Value = 172544863
$ws.Cells.Item(110, 3).Value = 34
$ws.Cells.Item(110, 4).Value = 44717
$ws.Cells.Item(115, 3).Value = 52251
$ws.Cells.Item(115, 4).Value = 76598088
$ws.Cells.Item(117, 3).Value = 26641
$ws.Cells.Item(117, 4).Value = 38596802
$ws.Cells.Item(118, 3).Value = 1302
$ws.Cells.Item(118, 4).Value = 1782491
$ws.Cells.Item(121, 3).Value = 2205
$ws.Cells.Item(121, 4).Value = 3096355
$ws.Cells.Item(123, 3).Value = 494878
$ws.Cells.Item(123, 4).Value = 652616063
$ws.Cells.Item(128, 3).Value = 1364
$ws.Cells.Item(128, 4).Value = 2022311
$ws.Cells.Item(130, 3).Value = 205116
$ws.Cells.Item(130, 4).Value = 301522707
$ws.Cells.Item(131, 3).Value = 390
$ws.Cells.Item(131, 4).Value = 581790
$ws.Cells.Item(133, 3).Value = 177301
$ws.Cells.Item(133, 4).Value = 257706162
$ws.Cells.Item(136, 3).Value = 2825
$ws.Cells.Item(136, 4).Value = 3970785
$ws.Cells.Item(138, 3).Value = 6192
$ws.Cells.Item(138, 4).Value = 8748866
$ws.Cells.Item(141, 3).Value = 43912
$ws.Cells.Item(141, 4).Value = 58625733
$ws.Cells.Item(147, 3).Value = 13916
$ws.Cells.Item(147, 4).Value = 20410330
$ws.Cells.Item(154, 3).Value = 17288
$ws.Cells.Item(154, 4).Value = 22842544
$ws.Cells.Item(158, 3).Value = 7066
$ws.Cells.Item(158, 4).Value = 10276203
$ws.Cells.Item(160, 3).Value = 4921
$ws.Cells.Item(160, 4).Value = 7081742
$ws.Cells.Item(163, 3).Value = 261
$ws.Cells.Item(163, 4).Value = 373774
$ws.Cells.Item(165, 3).Value = 15577
$ws.Cells.Item(165, 4).Value = 22604690
$ws.Cells.Item(166, 3).Value = 1755
$ws.Cells.Item(166, 4).Value = 2610230
$ws.Cells.Item(170, 3).Value = 79
$ws.Cells.Item(170, 4).Value = 118449
$ws.Cells.Item(171, 3).Value = 86665
$ws.Cells.Item(171, 4).Value = 108415256
$ws.Cells.Item(172, 3).Value = 32
$ws.Cells.Item(172, 4).Value = 37159
$ws.Cells.Item(178, 3).Value = 33570
$ws.Cells.Item(178, 4).Value = 49231427
$ws.Cells.Item(180, 3).Value = 12849
$ws.Cells.Item(180, 4).Value = 18563568
$ws.Cells.Item(182, 3).Value = 1240
$ws.Cells.Item(182, 4).Value = 1735896
$ws.Cells.Item(184, 3).Value = 1613
$ws.Cells.Item(184, 4).Value = 2266302
$ws.Cells.Item(186, 3).Value = 235667
$ws.Cells.Item(186, 4).Value = 292982325
$ws.Cells.Item(194, 3).Value = 85945
$ws.Cells.Item(194, 4).Value = 125984583
$ws.Cells.Item(197, 3).Value = 32672
$ws.Cells.Item(197, 4).Value = 47021204
$ws.Cells.Item(200, 3).Value = 5066
$ws.Cells.Item(200, 4).Value = 7216393
$ws.Cells.Item(203, 3).Value = 4772
$ws.Cells.Item(203, 4).Value = 6606033
$ws.Cells.Item(206, 3).Value = 260564
$ws.Cells.Item(206, 4).Value = 322506502
$ws.Cells.Item(208, 3).Value = 251
$ws.Cells.Item(208, 4).Value = 359087
$ws.Cells.Item(213, 3).Value = 611
$ws.Cells.Item(213, 4).Value = 889378
$ws.Cells.Item(215, 3).Value = 94373
$ws.Cells.Item(215, 4).Value = 138063569
$ws.Cells.Item(216, 3).Value = 87
$ws.Cells.Item(216, 4).Value = 129699
$ws.Cells.Item(218, 3).Value = 50840
$ws.Cells.Item(218, 4).Value = 73473772
$ws.Cells.Item(221, 3).Value = 4644
$ws.Cells.Item(221, 4).Value = 6520104
$ws.Cells.Item(224, 3).Value = 5613
$ws.Cells.Item(224, 4).Value = 7761659
$ws.Cells.Item(227, 3).Value = 104931
$ws.Cells.Item(227, 4).Value = 131294093
$ws.Cells.Item(234, 3).Value = 49094
$ws.Cells.Item(234, 4).Value = 71924456
$ws.Cells.Item(236, 3).Value = 12227
$ws.Cells.Item(236, 4).Value = 17578290
$ws.Cells.Item(238, 3).Value = 1883
$ws.Cells.Item(238, 4).Value = 2698784
$ws.Cells.Item(240, 3).Value = 2445
$ws.Cells.Item(240, 4).Value = 3416815
$ws.Cells.Item(241, 3).Value = 254108
$ws.Cells.Item(241, 4).Value = 320856273
$ws.Cells.Item(247, 3).Value = 820
$ws.Cells.Item(247, 4).Value = 1204550
$ws.Cells.Item(249, 3).Value = 94908
$ws.Cells.Item(249, 4).Value = 139069564
$ws.Cells.Item(252, 3).Value = 64129
$ws.Cells.Item(252, 4).Value = 92930689
$ws.Cells.Item(254, 3).Value = 2393
$ws.Cells.Item(254, 4).Value = 3376428
$ws.Cells.Item(257, 3).Value = 4510
$ws.Cells.Item(257, 4).Value = 6331539
